$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "OUT -> 2017/03/20 13:09",
    "IN -> 2017/03/20 13:11",
    "OUT -> 2017/03/20 13:11",
    "IN -> 2017/03/20 15:57",
    "OUT -> 2017/03/20 20:56",
    "IN -> 2017/03/20 20:56",
    "OUT -> 2017/03/20 20:57",
    "IN -> 2017/03/20 20:58",
    "OUT -> 2017/03/20 21:43",
    "IN -> 2017/03/20 21:43",
    "OUT -> 2017/03/20 21:43",
    "IN -> 2017/03/20 21:43",
    "OUT -> 2017/03/20 22:50",
    "IN -> 2017/03/20 22:50",
    "OUT -> 2017/03/20 22:50",
    "IN -> 2017/03/20 23:04"
)

# Make the newly added cells use the same (default) format as the
# existing data cells above them, rather than inheriting the column's
# default style, by copying the format from the last existing row.
$srcRange = $ws.Range("A41")
$destRange = $ws.Range("A42:A57")
$srcRange.Copy()
$destRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$row = 42
foreach ($val in $values) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $val
    $row = $row + 1
}
